$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2133333333333333
$ws.Range("C2").Value = 0.5233333333333333
$ws.Range("J2").Value = 0.02666666666666667
$ws.Range("P2").Value = 0.1433333333333333
$ws.Range("S2").Value = 0.09333333333333334

# Row 3
$ws.Range("B3").Value = 0.00625
$ws.Range("C3").Value = 0.03125
$ws.Range("J3").Value = 0.03125
$ws.Range("P3").Value = 0.69375
$ws.Range("S3").Value = 0.2375

# Row 4
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.6590909090909091
$ws.Range("S4").Value = 0.3181818181818182

# Row 6
$ws.Range("B6").Value = 0.02803738317757009
$ws.Range("D6").Value = 0.009345794392523364
$ws.Range("F6").Value = 0.06074766355140187
$ws.Range("J6").Value = 0.2663551401869159
$ws.Range("O6").Value = 0.04205607476635514
$ws.Range("Q6").Value = 0.1448598130841121
$ws.Range("R6").Value = 0.07943925233644859
$ws.Range("S6").Value = 0.3691588785046729

# Row 7
$ws.Range("B7").Value = 0.1302083333333333
$ws.Range("D7").Value = 0.02083333333333333
$ws.Range("F7").Value = 0.03125
$ws.Range("J7").Value = 0.1458333333333333
$ws.Range("O7").Value = 0.01041666666666667
$ws.Range("Q7").Value = 0.2083333333333333
$ws.Range("R7").Value = 0.08854166666666667
$ws.Range("S7").Value = 0.3645833333333333

# Row 8
$ws.Range("B8").Value = 0.08439897698209718
$ws.Range("D8").Value = 0.02557544757033248
$ws.Range("F8").Value = 0.05882352941176471
$ws.Range("J8").Value = 0.1534526854219949
$ws.Range("O8").Value = 0.01023017902813299
$ws.Range("Q8").Value = 0.1739130434782609
$ws.Range("R8").Value = 0.08184143222506395
$ws.Range("S8").Value = 0.4117647058823529

# Row 9
$ws.Range("B9").Value = 0.1041666666666667
$ws.Range("D9").Value = 0.02083333333333333
$ws.Range("F9").Value = 0.09375
$ws.Range("J9").Value = 0.1041666666666667
$ws.Range("O9").Value = 0.015625
$ws.Range("Q9").Value = 0.1614583333333333
$ws.Range("R9").Value = 0.078125
$ws.Range("S9").Value = 0.421875

# Row 10
$ws.Range("B10").Value = 0.1222130470685384
$ws.Range("D10").Value = 0.01981833195706028
$ws.Range("F10").Value = 0.07597027250206441
$ws.Range("J10").Value = 0.1098265895953757
$ws.Range("O10").Value = 0.0189925681255161
$ws.Range("Q10").Value = 0.1824938067712634
$ws.Range("R10").Value = 0.07018992568125516
$ws.Range("S10").Value = 0.4004954582989265

# Row 11
$ws.Range("G11").Value = 0.1561461794019934
$ws.Range("J11").Value = 0.106312292358804
$ws.Range("K11").Value = 0.2192691029900332
$ws.Range("L11").Value = 0.5116279069767442
$ws.Range("S11").Value = 0.006644518272425249

# Row 12
$ws.Range("G12").Value = 0.7098765432098766
$ws.Range("J12").Value = 0.1975308641975309
$ws.Range("K12").Value = 0.006172839506172839
$ws.Range("L12").Value = 0.04320987654320987
$ws.Range("S12").Value = 0.04320987654320987

# Row 13
$ws.Range("G13").Value = 0.7674418604651163
$ws.Range("J13").Value = 0.2325581395348837

# Row 14
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25

# Row 15
$ws.Range("F15").Value = 0.01932367149758454
$ws.Range("H15").Value = 0.1256038647342995
$ws.Range("I15").Value = 0.07729468599033816
$ws.Range("J15").Value = 0.3429951690821256
$ws.Range("K15").Value = 0.05314009661835749
$ws.Range("M15").Value = 0.03381642512077294
$ws.Range("O15").Value = 0.03381642512077294
$ws.Range("S15").Value = 0.3140096618357488

# Row 16
$ws.Range("F16").Value = 0.01657458563535912
$ws.Range("H16").Value = 0.1546961325966851
$ws.Range("I16").Value = 0.05524861878453038
$ws.Range("J16").Value = 0.4143646408839779
$ws.Range("K16").Value = 0.138121546961326
$ws.Range("O16").Value = 0.07734806629834254
$ws.Range("S16").Value = 0.143646408839779

# Row 17
$ws.Range("F17").Value = 0.02583979328165375
$ws.Range("H17").Value = 0.1679586563307494
$ws.Range("I17").Value = 0.07493540051679587
$ws.Range("J17").Value = 0.4366925064599483
$ws.Range("K17").Value = 0.09819121447028424
$ws.Range("M17").Value = 0.02583979328165375
$ws.Range("O17").Value = 0.06976744186046512
$ws.Range("S17").Value = 0.1007751937984496

# Row 18
$ws.Range("F18").Value = 0.01785714285714286
$ws.Range("H18").Value = 0.1785714285714286
$ws.Range("I18").Value = 0.08333333333333333
$ws.Range("J18").Value = 0.375
$ws.Range("K18").Value = 0.08928571428571429
$ws.Range("M18").Value = 0.01785714285714286
$ws.Range("O18").Value = 0.1011904761904762
$ws.Range("S18").Value = 0.1369047619047619

# Row 19
$ws.Range("F19").Value = 0.01609010458567981
$ws.Range("H19").Value = 0.1946902654867257
$ws.Range("I19").Value = 0.09814963797264682
$ws.Range("J19").Value = 0.3748994368463395
$ws.Range("K19").Value = 0.1150442477876106
$ws.Range("M19").Value = 0.01850362027353178
$ws.Range("N19").Value = 0.003218020917135961
$ws.Range("O19").Value = 0.06194690265486726
$ws.Range("S19").Value = 0.1174577634754626

